$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:G2").Value = 17.29
$ws.Range("B3:G3").Value = 9.83
